# Apply the ch-elm version 1.2.0 update to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": bump Version and Date values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.2.0"
$meta.Range("B8").Value = "2024-03-28T10:46:20+01:00"

# --- Sheet "Include from LOINC": insert new LOINC concept rows ---
$loinc = $wb.Worksheets.Item("Include from LOINC")

# Existing data (before edit):
# Row1: Concept | Description            (header, style s=1)
# Row2: 85827-4 | Carbapenem resistance bla OXA-48-like gene [Presence] by Molecular method  (style s=2)
# Row3: (empty) | (empty)                                                                    (style s=2)
# Row4: System URI | http://loinc.org                                                        (style s=2)

# First, extend the formatting (style s=2, same as existing data rows) down to the
# new rows 5-9 by copying the format of row 2.
$loinc.Range("A2:B2").Copy()
$loinc.Range("A5:B9").PasteSpecial(-4122)  # xlPasteFormats

# Now populate the values. The LOINC concepts are re-sorted and a new one is added,
# pushing the separator/footer rows further down.
$loinc.Range("A2").Value = "100911-7"
$loinc.Range("B2").Value = "Campylobacter sp [Presence] in Stool by Organism specific culture"

$loinc.Range("A3").Value = "49614-1"
$loinc.Range("B3").Value = "Campylobacter sp DNA [Identifier] in Specimen by NAA with probe detection"

$loinc.Range("A4").Value = "4992-4"
$loinc.Range("B4").Value = "Campylobacter sp rRNA [Presence] in Specimen by Probe"

$loinc.Range("A5").Value = "71429-5"
$loinc.Range("B5").Value = "Campylobacter sp DNA.diarrheagenic [Presence] in Stool by NAA with probe detection"

$loinc.Range("A6").Value = "85827-4"
$loinc.Range("B6").Value = "Carbapenem resistance bla OXA-48-like gene [Presence] by Molecular method"

$loinc.Range("A7").Value = "97513-6"
$loinc.Range("B7").Value = "Campylobacter sp [Presence] in Specimen by Organism specific culture"

$loinc.Range("A8").Value = ""
$loinc.Range("B8").Value = ""

$loinc.Range("A9").Value = "System URI"
$loinc.Range("B9").Value = "http://loinc.org"
